$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "train-exp" task names (row 8 edited before row 7)
$ws.Range("A8").Value = "宇宙黑洞"
$ws.Range("A7").Value = "时空相机"

# Add a new row 17 with an additional general-type task
# (type filled in before the game name)
$ws.Range("B17").Value = "general"
$ws.Range("A17").Value = "位置记忆PRO"

# Update the "test-near" task names
$ws.Range("A11").Value = "格子卡片"

# Update the "test-far" task names (row 14 edited before row 13)
$ws.Range("A14").Value = "数字卡片"
$ws.Range("A13").Value = "速算师（中级）"

# Restore the "test-near" task previously named differently (row 12)
$ws.Range("A12").Value = "幸运小球"

# Update the selected cell to match the new active cell
$ws.Range("D8").Select()
